# B6-PowerPoint.pptx edit
#
# The authored change (per the OOXML diff) does two things:
#
#   1. Re-colours the deck's theme so the slide master's theme
#      ("Integral" / "Red Violet" colour scheme) is replaced by the
#      standard Office colour scheme (the colours that, before this
#      edit, only lived in the unused ppt/theme/theme2.xml part used
#      by the notes master). Font scheme / format scheme are already
#      identical between the two theme parts, so only the 12 theme
#      colours need to change.
#
#   2. Switches the table style used by the three tables on slides
#      14-16 from the deck's custom "Table_0" style
#      ({F5E02AFD-4A60-45C7-8234-5F0BD6A9EFE5}) to the built-in
#      "No Style, Table Grid" style
#      ({E7703543-E7B1-4E94-88E9-C99E4F70A5A1}).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1. Recolour the presentation theme to the standard Office palette.
# ---------------------------------------------------------------

function Convert-HexToComRgb([string]$hex) {
    # PowerPoint's ColorFormat/RGBColor.RGB is an 0x00BBGGRR packed
    # integer (low byte = red), matching the classic VBA RGB() macro.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target values are the standard "Office" theme colours, in the
# ThemeColorScheme.Item(1..12) order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$theme = $p.Designs.Item(1).SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $colorScheme.Item($i).RGB = Convert-HexToComRgb $officeThemeHex[$i - 1]
}

# ---------------------------------------------------------------
# 2. Re-style every table that still uses the deck's custom
#    "Table_0" style to the built-in "No Style, Table Grid" style
#    (this lands on the three tables found on slides 14, 15 and 16).
# ---------------------------------------------------------------

$oldTableStyleId = "{F5E02AFD-4A60-45C7-8234-5F0BD6A9EFE5}"
$newTableStyleId = "{E7703543-E7B1-4E94-88E9-C99E4F70A5A1}"

for ($slideIndex = 1; $slideIndex -le $p.Slides.Count; $slideIndex++) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable -and $shape.Table.Style -eq $oldTableStyleId) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}
